$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.179.40"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "'2.058.08"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.46"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'16.19"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "'0.919"
$ws.Range("E13").Value = "  +13.89%  "
$ws.Range("D14").Value = "'2.358.63"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "'5.74"
$ws.Range("D16").Value = "'2.057.73"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "'18.94"
$ws.Range("E17").Value = "  +12.45%  "
$ws.Range("D18").Value = "'37.185.81"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'75.03"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'5.50"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "'238.21"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").Value = "'9.65"
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  -3.51%  "
$ws.Range("D27").Value = "'170.74"
$ws.Range("D28").Value = "'20.26"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'5.16"
$ws.Range("E30").Value = "  +8.92%  "
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'4.66"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D36").Value = "'2.31"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "'5.22"
$ws.Range("E39").Value = "  +11.41%  "
$ws.Range("E40").Value = "  -8.81%  "
$ws.Range("E41").Value = "  +7.85%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'17.61"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "'1.17"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").Value = "'96.99"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'2.43"
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "'1.278.88"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "'2.246.71"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'0.148"
$ws.Range("E51").Value = "  +9.93%  "
